# PROS-9738 - CCRU - New POS 2019 KPIs
# Inserts a new "Product Category" column before the existing "Brand" column (L),
# shifting L:Z -> M:AA, then fills in the new column's data for the two KPI rows
# that got a Product Category value, and updates the "Type" (CAT) values and the
# "Locations to include" text for KPI #5 (Juice Displays -> Juice Displays, Mixability Displays).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Target Execution 2019")

# Insert a new column before column L ("Brand"); this shifts Brand..level from L..Z to M..AA
$ws.Columns("L").Insert()

# Header for the new column
$ws.Range("L1").Value = "Product Category"

# Row 2 (KPI ID 1 - SSD Displays): Type changes from LOCATION_TYPE to CAT, new Product Category = SSD
$ws.Range("I2").Value = "CAT"
$ws.Range("L2").Value = "SSD"

# Row 6 (KPI ID 5 - Juice Displays): Type changes from LOCATION_TYPE to CAT, new Product Category = Juices,
# and the Locations to include value (now column P after the insert) gets a more specific label
$ws.Range("I6").Value = "CAT"
$ws.Range("L6").Value = "Juices"
$ws.Range("P6").Value = "Juice Displays, Mixability Displays"

# Highlight the new/changed data cells with the workbook's existing "yellow highlight" style
# (same formatting already used for similarly-flagged cells, e.g. the old O2/now P2).
$highlight = $ws.Range("P2").Interior.Color
foreach ($addr in @("I2","L2","I6","L6","P6")) {
    $ws.Range($addr).Interior.Color = $highlight
}

# Keep the AutoFilter in sync with the new used range
$ws.Range("A1:AA1").AutoFilter(1) | Out-Null
